$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_Details")

# Fill in the Execute_Flag values for the two test case rows
$ws.Range("C2").Value = "Yes"
$ws.Range("C3").Value = "No"

# Update the active selection to reflect the last edited cell
$ws.Activate()
$ws.Range("C3").Select()
